$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.717547
$ws.Range("H2").Value = 17.152641
$ws.Range("I2").Value = 0.3498537173739997
$ws.Range("J2").Value = 0.3498537173739997
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.248835333333334
$ws.Range("N2").Value = 6.746506
$ws.Range("O2").Value = 0.03590294220158827
$ws.Range("P2").Value = 0.03590294220158827
$ws.Range("Q2").Value = 12.857821713594
$ws.Range("R2").Value = 115.720395422346
$ws.Range("S2").Value = 0.01256077779388951
$ws.Range("T2").Value = 0.01256077779388951
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.717547
$ws.Range("H3").Value = 17.152641
$ws.Range("I3").Value = 0.3498537173739997
$ws.Range("J3").Value = 0.3498537173739997
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 44.29005966666667
$ws.Range("N3").Value = 132.870179
$ws.Range("O3").Value = 0.7070964373190639
$ws.Range("P3").Value = 0.7070964373190639
$ws.Range("Q3").Value = 253.230497776971
$ws.Range("R3").Value = 2279.074479992739
$ws.Range("S3").Value = 0.2473803171379859
$ws.Range("T3").Value = 0.2473803171379859
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.717547
$ws.Range("H4").Value = 17.152641
$ws.Range("I4").Value = 0.3498537173739997
$ws.Range("J4").Value = 0.3498537173739997
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.09762433333333
$ws.Range("N4").Value = 48.292873
$ws.Range("O4").Value = 0.2570006204793478
$ws.Range("P4").Value = 0.2570006204793479
$ws.Range("Q4").Value = 92.03892371417699
$ws.Range("R4").Value = 828.3503134275929
$ws.Range("S4").Value = 0.08991262244212429
$ws.Range("T4").Value = 0.08991262244212432
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.260042333333333
$ws.Range("H5").Value = 9.780127
$ws.Range("I5").Value = 0.1994802892067655
$ws.Range("J5").Value = 0.1994802892067655
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.248835333333334
$ws.Range("N5").Value = 6.746506
$ws.Range("O5").Value = 0.03590294220158827
$ws.Range("P5").Value = 0.03590294220158827
$ws.Range("Q5").Value = 7.331298387362445
$ws.Range("R5").Value = 65.98168548626201
$ws.Range("S5").Value = 0.007161929293746615
$ws.Range("T5").Value = 0.007161929293746616
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.260042333333333
$ws.Range("H6").Value = 9.780127
$ws.Range("I6").Value = 0.1994802892067655
$ws.Range("J6").Value = 0.1994802892067655
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.29005966666667
$ws.Range("N6").Value = 132.870179
$ws.Range("O6").Value = 0.7070964373190639
$ws.Range("P6").Value = 0.7070964373190639
$ws.Range("Q6").Value = 144.3874694591926
$ws.Range("R6").Value = 1299.487225132733
$ws.Range("S6").Value = 0.1410518018134804
$ws.Range("T6").Value = 0.1410518018134804
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.260042333333333
$ws.Range("H7").Value = 9.780127
$ws.Range("I7").Value = 0.1994802892067655
$ws.Range("J7").Value = 0.1994802892067655
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.09762433333333
$ws.Range("N7").Value = 48.292873
$ws.Range("O7").Value = 0.2570006204793478
$ws.Range("P7").Value = 0.2570006204793479
$ws.Range("Q7").Value = 52.47893679276344
$ws.Range("R7").Value = 472.310431134871
$ws.Range("S7").Value = 0.05126655809953848
$ws.Range("T7").Value = 0.0512665580995385
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.093744666666666
$ws.Range("H8").Value = 6.281234
$ws.Range("I8").Value = 0.1281151435861076
$ws.Range("J8").Value = 0.1281151435861077
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.248835333333334
$ws.Range("N8").Value = 6.746506
$ws.Range("O8").Value = 0.03590294220158827
$ws.Range("P8").Value = 0.03590294220158827
$ws.Range("Q8").Value = 4.708486985378222
$ws.Range("R8").Value = 42.376382868404
$ws.Range("S8").Value = 0.004599710595320205
$ws.Range("T8").Value = 0.004599710595320206
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.093744666666666
$ws.Range("H9").Value = 6.281234
$ws.Range("I9").Value = 0.1281151435861076
$ws.Range("J9").Value = 0.1281151435861077
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 44.29005966666667
$ws.Range("N9").Value = 132.870179
$ws.Range("O9").Value = 0.7070964373190639
$ws.Range("P9").Value = 0.7070964373190639
$ws.Range("Q9").Value = 92.73207621343177
$ws.Range("R9").Value = 834.588685920886
$ws.Range("S9").Value = 0.09058976159635704
$ws.Range("T9").Value = 0.09058976159635707
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.093744666666666
$ws.Range("H10").Value = 6.281234
$ws.Range("I10").Value = 0.1281151435861076
$ws.Range("J10").Value = 0.1281151435861077
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.09762433333333
$ws.Range("N10").Value = 48.292873
$ws.Range("O10").Value = 0.2570006204793478
$ws.Range("P10").Value = 0.2570006204793479
$ws.Range("Q10").Value = 33.70431509392021
$ws.Range("R10").Value = 303.338835845282
$ws.Range("S10").Value = 0.0329256713944304
$ws.Range("T10").Value = 0.03292567139443042
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.630401
$ws.Range("H11").Value = 7.891203
$ws.Range("I11").Value = 0.1609528645823613
$ws.Range("J11").Value = 0.1609528645823613
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.248835333333334
$ws.Range("N11").Value = 6.746506
$ws.Range("O11").Value = 0.03590294220158827
$ws.Range("P11").Value = 0.03590294220158827
$ws.Range("Q11").Value = 5.915338709635334
$ws.Range("R11").Value = 53.238048386718
$ws.Range("S11").Value = 0.005778681394280582
$ws.Range("T11").Value = 0.005778681394280583
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.630401
$ws.Range("H12").Value = 7.891203
$ws.Range("I12").Value = 0.1609528645823613
$ws.Range("J12").Value = 0.1609528645823613
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.29005966666667
$ws.Range("N12").Value = 132.870179
$ws.Range("O12").Value = 0.7070964373190639
$ws.Range("P12").Value = 0.7070964373190639
$ws.Range("Q12").Value = 116.5006172372597
$ws.Range("R12").Value = 1048.505555135337
$ws.Range("S12").Value = 0.1138091971224854
$ws.Range("T12").Value = 0.1138091971224854
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.630401
$ws.Range("H13").Value = 7.891203
$ws.Range("I13").Value = 0.1609528645823613
$ws.Range("J13").Value = 0.1609528645823613
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.09762433333333
$ws.Range("N13").Value = 48.292873
$ws.Range("O13").Value = 0.2570006204793478
$ws.Range("P13").Value = 0.2570006204793479
$ws.Range("Q13").Value = 42.34320714402433
$ws.Range("R13").Value = 381.088864296219
$ws.Range("S13").Value = 0.04136498606559529
$ws.Range("T13").Value = 0.04136498606559531
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.640944
$ws.Range("H14").Value = 7.922832
$ws.Range("I14").Value = 0.1615979852507658
$ws.Range("J14").Value = 0.1615979852507658
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.248835333333334
$ws.Range("N14").Value = 6.746506
$ws.Range("O14").Value = 0.03590294220158827
$ws.Range("P14").Value = 0.03590294220158827
$ws.Range("Q14").Value = 5.939048180554667
$ws.Range("R14").Value = 53.451433624992
$ws.Range("S14").Value = 0.005801843124351357
$ws.Range("T14").Value = 0.005801843124351359
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.640944
$ws.Range("H15").Value = 7.922832
$ws.Range("I15").Value = 0.1615979852507658
$ws.Range("J15").Value = 0.1615979852507658
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 44.29005966666667
$ws.Range("N15").Value = 132.870179
$ws.Range("O15").Value = 0.7070964373190639
$ws.Range("P15").Value = 0.7070964373190639
$ws.Range("Q15").Value = 116.9675673363253
$ws.Range("R15").Value = 1052.708106026928
$ws.Range("S15").Value = 0.1142653596487551
$ws.Range("T15").Value = 0.1142653596487551
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.640944
$ws.Range("H16").Value = 7.922832
$ws.Range("I16").Value = 0.1615979852507658
$ws.Range("J16").Value = 0.1615979852507658
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.09762433333333
$ws.Range("N16").Value = 48.292873
$ws.Range("O16").Value = 0.2570006204793478
$ws.Range("P16").Value = 0.2570006204793479
$ws.Range("Q16").Value = 33.70431509392021
$ws.Range("R16").Value = 303.338835845282
$ws.Range("S16").Value = 0.0415307824776593
$ws.Range("T16").Value = 0.0415307824776593
